{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Start after the last existing paragraph (the \"[5] ...\" source paragraph).\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nconst sources = [\n  {\n    label: \"[6] \",\n    url: \"https://www.unsw.edu.au/research/ndarc/news-events/blogs/2015/05/don_t-believe-the-hype--teens-are-drinking-less-than-they-used-t\",\n    splitLabel: false\n  },\n  {\n    label: \"[7]\",\n    url: \"https://www.tobaccoinaustralia.org.au/chapter-18-e-cigarettes/18-12-public-perceptions-of-the-risks-and-benefits-of-e-cigarettes\",\n    splitLabel: true\n  },\n  {\n    label: \"[8]\",\n    url: \"https://www.cancer.nsw.gov.au/prevention-and-screening/preventing-cancer/damaging-effects-of-vaping/vaping-harms-your-health\",\n    splitLabel: true\n  },\n  {\n    label: \"[9] \",\n    url: \"https://www.health.gov.au/news/new-national-campaign-launched-to-help-young-people-quit-vaping-0\",\n    splitLabel: false\n  }\n];\n\nfor (const source of sources) {\n  const newPara = anchor.insertParagraph(\"\", Word.InsertLocation.after);\n\n  if (source.splitLabel) {\n    // e.g. \"[7]\" then a separate \" \" run, mirroring the source document.\n    newPara.insertText(source.label, Word.InsertLocation.start);\n    newPara.insertText(\" \", Word.InsertLocation.end);\n  } else {\n    // e.g. \"[6] \" as a single run with trailing space included.\n    newPara.insertText(source.label, Word.InsertLocation.start);\n  }\n\n  const linkRange = newPara.insertText(source.url, Word.InsertLocation.end);\n  linkRange.hyperlink = source.url;\n\n  newPara.insertText(\" \", Word.InsertLocation.end);\n\n  anchor = newPara;\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$sources = @(\n  @{ Label = \"[6] \"; Url = \"https://www.unsw.edu.au/research/ndarc/news-events/blogs/2015/05/don_t-believe-the-hype--teens-are-drinking-less-than-they-used-t\" },\n  @{ Label = \"[7] \"; Url = \"https://www.tobaccoinaustralia.org.au/chapter-18-e-cigarettes/18-12-public-perceptions-of-the-risks-and-benefits-of-e-cigarettes\" },\n  @{ Label = \"[8] \"; Url = \"https://www.cancer.nsw.gov.au/prevention-and-screening/preventing-cancer/damaging-effects-of-vaping/vaping-harms-your-health\" },\n  @{ Label = \"[9] \"; Url = \"https://www.health.gov.au/news/new-national-campaign-launched-to-help-young-people-quit-vaping-0\" }\n)\n\nforeach ($source in $sources) {\n  # Start a new paragraph right after the current end of the document.\n  $endPos = $d.Content.End\n  $r = $d.Range($endPos, $endPos)\n  $r.InsertParagraphAfter()\n\n  # Write the \"[n] \" label into the freshly created (empty) paragraph.\n  $labelPos = $d.Content.End - 1\n  $rLabel = $d.Range($labelPos, $labelPos)\n  $rLabel.InsertAfter($source.Label)\n\n  # Insert the URL text, then convert that exact span into a hyperlink.\n  $urlStart = $d.Content.End - 1\n  $rUrl = $d.Range($urlStart, $urlStart)\n  $rUrl.InsertAfter($source.Url)\n  $urlEnd = $urlStart + $source.Url.Length\n  $linkRange = $d.Range($urlStart, $urlEnd)\n  $d.Hyperlinks.Add($linkRange, $source.Url, \"\", \"\", $source.Url) | Out-Null\n\n  # Trailing space after the hyperlink, still inside the paragraph.\n  $trailPos = $d.Content.End - 1\n  $rTrail = $d.Range($trailPos, $trailPos)\n  $rTrail.InsertAfter(\" \")\n}\n"}
